# S_PM_05.xlsx edit script
# Implements the changes described by the commit:
#  - Developing&Training_Cost (sheet1): update hours/cost figures, apply a
#    thousands-separator number format to the Cost($) column, and change
#    the selection/view state.
#  - Tools_Cost (sheet2): bump the SW tools cost, apply the same number
#    format, and change the selection.
#  - Total_Cost (sheet3): bump the overall total, apply the same number
#    format, and change the selection.
#  - Assumptions (sheet4): no data changes, but it is no longer the
#    active/focused sheet.

$wb = $excel.ActiveWorkbook

$wsDev  = $wb.Worksheets.Item("Developing&Training_Cost")
$wsTool = $wb.Worksheets.Item("Tools_Cost")
$wsTot  = $wb.Worksheets.Item("Total_Cost")
$wsAsm  = $wb.Worksheets.Item("Assumptions")

# ---- Developing&Training_Cost ----------------------------------------
$wsDev.Range("J6").Value = 480
$wsDev.Range("K6").Value = 28800

$wsDev.Range("J7").Value = 160
$wsDev.Range("K7").Value = 8000

$wsDev.Range("J8").Value = 80
$wsDev.Range("K8").Value = 3200

$wsDev.Range("K6:K8").NumberFormat = "#,##0"

# ---- Tools_Cost --------------------------------------------------------
$wsTool.Range("B7").Value = 10000
$wsTool.Range("B7").NumberFormat = "#,##0"

# ---- Total_Cost ---------------------------------------------------------
$wsTot.Range("A2").Value = 50000
$wsTot.Range("A2").NumberFormat = "#,##0"

# ---- view / selection state --------------------------------------------
$wsTool.Range("B7").Select()
$wsTot.Range("A2").Select()
$wsAsm.Range("C5").Select()

# Developing&Training_Cost becomes the active sheet, with K8 selected.
$wsDev.Activate()
$wsDev.Range("K8").Select()
